# Apply updated cryptocurrency market data (prices and 1h volume change)
# generated from the authoritative diff of xl/worksheets/sheet1.xml
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.474.51'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '2.422.94'
$ws.Range("E3").Value = '  +2.97%  '
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '''510.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '''133.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.76%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = '2.449.96'
$ws.Range("E9").Value = '  +3.50%  '
$ws.Range("D10").Value = '''0.0973'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.47%  '
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("D12").Value = '''0.324'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.82%  '
$ws.Range("D13").Value = '''4.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.28%  '
$ws.Range("D14").Value = '2.881.83'
$ws.Range("E14").Value = '  +4.06%  '
$ws.Range("D15").Value = '57.412.54'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("E16").Value = '  +3.06%  '
$ws.Range("D17").Value = '''0.0000135'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.91%  '
$ws.Range("D18").Value = '2.491.64'
$ws.Range("E18").Value = '  +5.02%  '
$ws.Range("D19").Value = '''10.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").Value = '''313.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").Value = '''4.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("E22").Value = '  +5.36%  '
$ws.Range("D23").Value = '''5.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").Value = '''65.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("D26").Value = '''0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = '''0.383'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.155'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("E29").Value = '  +6.17%  '
$ws.Range("D30").Value = '''170.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("E31").Value = '  +2.79%  '
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("D33").Value = '''6.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("E34").Value = '  +0.27%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").Value = '''0.994'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").Value = '''18.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.12%  '
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("D39").Value = '''3.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.74%  '
$ws.Range("D40").Value = '''36.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.48%  '
$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D41").Value = '''0.805'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.60%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.45'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("D43").Value = '''132.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.03%  '
$ws.Range("E44").Value = '  +2.83%  '
$ws.Range("D45").Value = '''4.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").Value = '''254.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("D47").Value = '''0.570'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").Value = '''0.0915'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("D50").Value = '''17.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.54%  '
$ws.Range("E51").Value = '  +3.00%  '
